$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 15
$ws.Range("B1").Value = 3.010123014450073
$ws.Range("C1").Value = 2.665281057357788
$ws.Range("D1").Value = 2.936100006103516
$ws.Range("E1").Value = 15
